$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-19 Monday" "2025-05-20 Tuesday"

Replace-Text "256×6=1536" "464×6=2784"
Replace-Text "937×4=3748" "619×9=5571"
Replace-Text "666×8=5328" "913×5=4565"
Replace-Text "681×5=3405" "555×3=1665"
Replace-Text "532×9=4788" "589×3=1767"

Replace-Text "852×2=1704" "243×4=972"
Replace-Text "521×2=1042" "280×7=1960"
Replace-Text "221×9=1989" "281×6=1686"
Replace-Text "469×6=2814" "288×4=1152"
Replace-Text "897×4=3588" "853×9=7677"

Replace-Text "276×4=1104" "391×5=1955"
Replace-Text "219×6=1314" "886×7=6202"
Replace-Text "498×6=2988" "194×7=1358"
Replace-Text "190×7=1330" "817×5=4085"
Replace-Text "888×8=7104" "353×7=2471"

Replace-Text "669×4=2676" "186×6=1116"
Replace-Text "641×9=5769" "558×9=5022"
Replace-Text "251×9=2259" "706×6=4236"
Replace-Text "709×5=3545" "179×2=358"
Replace-Text "650×6=3900" "485×6=2910"

Replace-Text "644×6=3864" "331×8=2648"
Replace-Text "218×9=1962" "611×3=1833"
Replace-Text "432×5=2160" "168×7=1176"
Replace-Text "236×2=472" "386×7=2702"
Replace-Text "992×8=7936" "568×6=3408"
